# Auto-generated edit script applying numeric updates per the commit diff.
# Each statement sets a cell's Value on the appropriate worksheet.
# Cells whose new value is blank/removed in the diff are cleared (set to $null).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 13166.5
$ws.Range("J3").Value = 13166.5
$ws.Range("L3").Value = 13166.5
$ws.Range("N3").Value = -13394.5
$ws.Range("H6").Value = 1128
$ws.Range("I6").Value = 1128
$ws.Range("K6").Value = 3384
$ws.Range("M6").Value = -3272
$ws.Range("H12").Value = 8712307
$ws.Range("I12").Value = 14520421
$ws.Range("K12").Value = 14520421
$ws.Range("M12").Value = -14520251
$ws.Range("H17").Value = 3054303.2
$ws.Range("J17").Value = 3187055.5
$ws.Range("L17").Value = 9561166.5
$ws.Range("N17").Value = -9561502.5
$ws.Range("H28").Value = 398
$ws.Range("I28").Value = 452.5
$ws.Range("K28").Value = 452.5
$ws.Range("M28").Value = 32.5
$ws.Range("H52").Value = 2371.1428
$ws.Range("J52").Value = 1349.75
$ws.Range("L52").Value = 4049.25
$ws.Range("N52").Value = -4369.25
$ws.Range("H62").Value = 2829.6667
$ws.Range("I62").Value = 2632.3635
$ws.Range("K62").Value = 2632.3635
$ws.Range("M62").Value = -2008.3635
$ws.Range("H65").Value = 2829.6667
$ws.Range("I65").Value = 2632.3635
$ws.Range("K65").Value = 13161.8175
$ws.Range("M65").Value = -10041.8175
$ws.Range("H80").Value = 46297560
$ws.Range("J80").Value = 16668204
$ws.Range("L80").Value = 50004612
$ws.Range("N80").Value = -50006608
$ws.Range("H83").Value = 46297560
$ws.Range("J83").Value = 16668204
$ws.Range("L83").Value = 150013836
$ws.Range("N83").Value = -150023820
$ws.Range("H88").Value = 1829
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 1993.5
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 1993.5
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -2805.5
$ws.Range("H91").Value = 1829
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 1993.5
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 1993.5
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -4801.5
$ws.Range("H95").Value = 57499
$ws.Range("J95").Value = 57499
$ws.Range("L95").Value = 57499
$ws.Range("N95").Value = -62991
$ws.Range("H100").Value = 2002.2667
$ws.Range("J100").Value = 1598
$ws.Range("L100").Value = 1598
$ws.Range("N100").Value = -2680
$ws.Range("H102").Value = 13166.5
$ws.Range("J102").Value = 13166.5
$ws.Range("L102").Value = 13166.5
$ws.Range("N102").Value = -19656.5
$ws.Range("H105").Value = 33408.832
$ws.Range("J105").Value = 33408.832
$ws.Range("L105").Value = 33408.832
$ws.Range("N105").Value = -40396.832
$ws.Range("H106").Value = 66669916
$ws.Range("I106").Value = 71431520
$ws.Range("K106").Value = 71431520
$ws.Range("M106").Value = -71430889
$ws.Range("H116").Value = 7834.4
$ws.Range("I116").Value = 7449.8335
$ws.Range("J116").Value = 8411.25
$ws.Range("K116").Value = 7449.8335
$ws.Range("L116").Value = 8411.25
$ws.Range("M116").Value = -4007.8335
$ws.Range("N116").Value = -15295.25
$ws.Range("H137").Value = 16668233
$ws.Range("I137").Value = 897
$ws.Range("J137").Value = 20001700
$ws.Range("K137").Value = 2691
$ws.Range("L137").Value = 60005100
$ws.Range("M137").Value = -141
$ws.Range("N137").Value = -60010200
$ws.Range("H138").Value = 5445.4136
$ws.Range("J138").Value = 3875.383
$ws.Range("L138").Value = 11626.149
$ws.Range("N138").Value = -21906.149
$ws.Range("H141").Value = 2704.1333
$ws.Range("I141").Value = 2704.1333
$ws.Range("K141").Value = 8112.3999
$ws.Range("M141").Value = -2932.3999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 12510596
$ws.Range("I29").Value = 20006054
$ws.Range("K29").Value = 20006054
$ws.Range("M29").Value = -20005746
$ws.Range("H32").Value = 189899.19
$ws.Range("I32").Value = 196826.33
$ws.Range("K32").Value = 196826.33
$ws.Range("M32").Value = -196539.33
$ws.Range("H33").Value = 33340000
$ws.Range("I33").Value = 33340000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 33340000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -33339671
$ws.Range("N33").Value = $null
$ws.Range("H61").Value = 1254838.4
$ws.Range("I61").Value = 3510.2122
$ws.Range("K61").Value = 3510.2122
$ws.Range("M61").Value = -3298.2122
$ws.Range("H74").Value = 892785.6
$ws.Range("I74").Value = 2214.889
$ws.Range("J74").Value = 1393731.8
$ws.Range("K74").Value = 2214.889
$ws.Range("L74").Value = 1393731.8
$ws.Range("M74").Value = -1340.889
$ws.Range("N74").Value = -1395479.8
$ws.Range("H77").Value = 892785.6
$ws.Range("I77").Value = 2214.889
$ws.Range("J77").Value = 1393731.8
$ws.Range("K77").Value = 11074.445
$ws.Range("L77").Value = 6968659
$ws.Range("M77").Value = -6706.445
$ws.Range("N77").Value = -6977395
$ws.Range("H95").Value = 3997
$ws.Range("J95").Value = 3997
$ws.Range("L95").Value = 3997
$ws.Range("N95").Value = -9489
$ws.Range("H104").Value = 42500
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 42500
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 42500
$ws.Range("M104").Value = $null
$ws.Range("N104").Value = -49488
$ws.Range("H132").Value = 1579.0952
$ws.Range("I132").Value = 1191.898
$ws.Range("K132").Value = 3575.694
$ws.Range("M132").Value = -1045.694
$ws.Range("H136").Value = 1254838.4
$ws.Range("I136").Value = 3510.2122
$ws.Range("K136").Value = 10530.6366
$ws.Range("M136").Value = -7980.6366

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = $null
$ws.Range("H103").Value = 15024.375
$ws.Range("J103").Value = 15024.375
$ws.Range("L103").Value = 15024.375
$ws.Range("N103").Value = -17368.375
$ws.Range("H105").Value = 9084.625
$ws.Range("I105").Value = 10746
$ws.Range("J105").Value = 5429.6
$ws.Range("K105").Value = 10746
$ws.Range("L105").Value = 5429.6
$ws.Range("M105").Value = -8999
$ws.Range("N105").Value = -8923.6
$ws.Range("H134").Value = 20932306
$ws.Range("I134").Value = 2028.6774
$ws.Range("K134").Value = 6086.0322
$ws.Range("M134").Value = -3551.0322

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 6750
$ws.Range("J8").Value = 6750
$ws.Range("L8").Value = 6750
$ws.Range("N8").Value = -7030
$ws.Range("H19").Value = 1998
$ws.Range("I19").Value = 372.75
$ws.Range("K19").Value = 372.75
$ws.Range("M19").Value = -202.75
$ws.Range("H24").Value = 1998
$ws.Range("I24").Value = 372.75
$ws.Range("K24").Value = 372.75
$ws.Range("M24").Value = -202.75
$ws.Range("H43").Value = 22608
$ws.Range("J43").Value = 22608
$ws.Range("L43").Value = 22608
$ws.Range("N43").Value = -22976
$ws.Range("H62").Value = 1855.7
$ws.Range("I62").Value = 1855.7
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1855.7
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1231.7
$ws.Range("N62").Value = $null
$ws.Range("H65").Value = 1855.7
$ws.Range("I65").Value = 1855.7
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 9278.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -6158.5
$ws.Range("N65").Value = $null
$ws.Range("H74").Value = 44847.668
$ws.Range("J74").Value = 44847.668
$ws.Range("L74").Value = 44847.668
$ws.Range("N74").Value = -46595.668
$ws.Range("H77").Value = 44847.668
$ws.Range("J77").Value = 44847.668
$ws.Range("L77").Value = 134543.004
$ws.Range("N77").Value = -143279.004
$ws.Range("H101").Value = 22608
$ws.Range("J101").Value = 22608
$ws.Range("L101").Value = 22608
$ws.Range("N101").Value = -29098
$ws.Range("H107").Value = 697.1739
$ws.Range("I107").Value = 837.5294
$ws.Range("K107").Value = 837.5294
$ws.Range("M107").Value = 1082.4706

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4889398
$ws.Range("I4").Value = 7167122.5
$ws.Range("J4").Value = 333949
$ws.Range("K4").Value = 21501367.5
$ws.Range("L4").Value = 1001847
$ws.Range("M4").Value = -21501255.5
$ws.Range("N4").Value = -1002071
$ws.Range("H131").Value = 4547779
$ws.Range("I131").Value = 8265907
$ws.Range("J131").Value = 3400
$ws.Range("K131").Value = 24797721
$ws.Range("L131").Value = 10200
$ws.Range("M131").Value = -24792681
$ws.Range("N131").Value = -20280
$ws.Range("H132").Value = 1822.1428
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 1959.1666
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 17632.4994
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -22692.4994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 9300856
$ws.Range("I11").Value = 10847667
$ws.Range("J11").Value = 19995
$ws.Range("K11").Value = 10847667
$ws.Range("L11").Value = 19995
$ws.Range("M11").Value = -10847528
$ws.Range("N11").Value = -20273
$ws.Range("H12").Value = 9000
$ws.Range("J12").Value = 9000
$ws.Range("L12").Value = 9000
$ws.Range("N12").Value = -9280
$ws.Range("H14").Value = 2100
$ws.Range("J14").Value = 2000
$ws.Range("L14").Value = 2000
$ws.Range("N14").Value = -2336
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").Value = $null
$ws.Range("H58").Value = 265000
$ws.Range("I58").Value = 30000
$ws.Range("J58").Value = 500000
$ws.Range("K58").Value = 30000
$ws.Range("L58").Value = 500000
$ws.Range("M58").Value = -29723
$ws.Range("N58").Value = -500554
$ws.Range("H80").Value = 10992721
$ws.Range("I80").Value = 116407.24
$ws.Range("K80").Value = 116407.24
$ws.Range("M80").Value = -115409.24
$ws.Range("H83").Value = 10992721
$ws.Range("I83").Value = 116407.24
$ws.Range("K83").Value = 582036.2000000001
$ws.Range("M83").Value = -577044.2000000001
$ws.Range("H101").Value = 33899.6
$ws.Range("J101").Value = 33899.6
$ws.Range("L101").Value = 33899.6
$ws.Range("N101").Value = -40389.6
$ws.Range("H126").Value = 2727.7693
$ws.Range("I126").Value = 2678.7273
$ws.Range("J126").Value = 2997.5
$ws.Range("K126").Value = 8036.1819
$ws.Range("L126").Value = 8992.5
$ws.Range("M126").Value = -5566.1819
$ws.Range("N126").Value = -13932.5
$ws.Range("H132").Value = 613089.5
$ws.Range("I132").Value = 2112.6155
$ws.Range("J132").Value = 1605927
$ws.Range("K132").Value = 6337.8465
$ws.Range("L132").Value = 4817781
$ws.Range("M132").Value = -3807.8465
$ws.Range("N132").Value = -4822841

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 10000
$ws.Range("J3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("N3").Value = -10224
$ws.Range("H15").Value = 10000
$ws.Range("J15").Value = 10000
$ws.Range("L15").Value = 10000
$ws.Range("N15").Value = -10340
$ws.Range("H17").Value = 8152
$ws.Range("I17").Value = 8202.666999999999
$ws.Range("K17").Value = 8202.666999999999
$ws.Range("M17").Value = -8032.666999999999
$ws.Range("H32").Value = 131
$ws.Range("I32").Value = 131
$ws.Range("K32").Value = 131
$ws.Range("M32").Value = 186
$ws.Range("H40").Value = 3524.4443
$ws.Range("I40").Value = 2113.7273
$ws.Range("K40").Value = 2113.7273
$ws.Range("M40").Value = -1977.7273
$ws.Range("H61").Value = 5028.4375
$ws.Range("I61").Value = 5885.6665
$ws.Range("J61").Value = 3926.2856
$ws.Range("K61").Value = 5885.6665
$ws.Range("L61").Value = 3926.2856
$ws.Range("M61").Value = -5683.6665
$ws.Range("N61").Value = -4330.2856
$ws.Range("H93").Value = 1868.4193
$ws.Range("I93").Value = 1384.9524
$ws.Range("J93").Value = 2883.7
$ws.Range("K93").Value = 1384.9524
$ws.Range("L93").Value = 2883.7
$ws.Range("M93").Value = -136.9523999999999
$ws.Range("N93").Value = -5379.7
$ws.Range("H103").Value = 31749.75
$ws.Range("J103").Value = 31749.75
$ws.Range("L103").Value = 31749.75
$ws.Range("N103").Value = -34093.75
$ws.Range("H113").Value = 5028.4375
$ws.Range("I113").Value = 5885.6665
$ws.Range("J113").Value = 3926.2856
$ws.Range("K113").Value = 5885.6665
$ws.Range("L113").Value = 3926.2856
$ws.Range("M113").Value = -3715.6665
$ws.Range("N113").Value = -8266.285599999999
$ws.Range("H136").Value = 3368.2046
$ws.Range("I136").Value = 2086
$ws.Range("K136").Value = 6258
$ws.Range("M136").Value = -3708

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 2500
$ws.Range("I9").Value = 1500
$ws.Range("J9").Value = 3500
$ws.Range("K9").Value = 1500
$ws.Range("L9").Value = 3500
$ws.Range("M9").Value = -1360
$ws.Range("N9").Value = -3780
$ws.Range("H81").Value = 2565.2307
$ws.Range("J81").Value = 1334
$ws.Range("L81").Value = 2668
$ws.Range("N81").Value = -4790
$ws.Range("H84").Value = 2565.2307
$ws.Range("J84").Value = 1334
$ws.Range("L84").Value = 13340
$ws.Range("N84").Value = -23948
$ws.Range("H92").Value = 75000
$ws.Range("J92").Value = 75000
$ws.Range("L92").Value = 75000
$ws.Range("N92").Value = -79992
$ws.Range("H97").Value = 25895.4
$ws.Range("J97").Value = 25895.4
$ws.Range("L97").Value = 25895.4
$ws.Range("N97").Value = -27877.4
$ws.Range("H98").Value = 60000
$ws.Range("J98").Value = 60000
$ws.Range("L98").Value = 60000
$ws.Range("N98").Value = -65990
$ws.Range("H104").Value = 37233
$ws.Range("J104").Value = 37233
$ws.Range("L104").Value = 37233
$ws.Range("N104").Value = -44221
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = $null
$ws.Range("H122").Value = 1487.0465
$ws.Range("I122").Value = 1375.7435
$ws.Range("K122").Value = 4127.2305
$ws.Range("M122").Value = -1677.2305
